$wb = $excel.ActiveWorkbook

# Sheets that contain the event data rows to be updated: 展览 (Exhibition) and 全部类型 (All types)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1128
    $ws.Range("F3").Value = 73
}
